$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RERC")

# Delete column C (the "Inicio de Operaciones" date column). This shifts
# column D ("Electricidad generada") into C and column E (the GHG emissions
# formula/total) into D.
$ws.Range("C:C").Delete()

# Leave the cursor on the cell to the right of the remaining table, matching
# where the author ended up after trimming the column.
$ws.Range("F8").Select() | Out-Null
